# Corrected some selection scopes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2; B = 0.08692415871688995;   C = 0.4489496882726305;  D = 0.3262969795491077;  E = 0.5712241062394932; F = 0.577260257759593;   G = 23 },
    @{ Row = 3; B = -0.002553697066920395; C = 0.5658358590261171;  D = 0.5376784437804385;  E = 0.7332656024800553; F = 0.7505167201924804;  G = 22 },
    @{ Row = 4; B = 0.0651984737114082;    C = 0.3928637431660651;  D = 0.2463207646819407;  E = 0.4963071273737067; F = 0.5041561463406211;  G = 21 },
    @{ Row = 5; B = 0.08449645211527543;   C = 0.5227040812392598;  D = 0.4080452692200807;  E = 0.638784211780536;  F = 0.6496198230706389;  G = 20 },
    @{ Row = 6; B = 0.008651778582802071;  C = 0.3534920523821091;  D = 0.1987829789312101;  E = 0.4458508483015482; F = 0.4579819493964025;  G = 19 },
    @{ Row = 7; B = 0.1021562244677468;    C = 0.5059495565191923;  D = 0.3884286911590425;  E = 0.6232404761879979; F = 0.6326354932905994;  G = 18 },
    @{ Row = 8; B = -0.01064431491224962;  C = 0.3832088633390937;  D = 0.1943466901180776;  E = 0.4408476949220418; F = 0.4542829244760465;  G = 17 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Range("B$r").Value = $entry.B
    $ws.Range("C$r").Value = $entry.C
    $ws.Range("D$r").Value = $entry.D
    $ws.Range("E$r").Value = $entry.E
    $ws.Range("F$r").Value = $entry.F
    $ws.Range("G$r").Value = $entry.G
}
